$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting all existing data rows (27-107) down to (28-108).
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new observation data.
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value = 44980
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 100112040
$ws.Cells.Item(27, 7).Value = "Cilantro"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 350
$ws.Cells.Item(27, 11).Value = 3500
$ws.Cells.Item(27, 12).Value = 4500
$ws.Cells.Item(27, 13).Value = 3929
$ws.Cells.Item(27, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 1964
$ws.Cells.Item(27, 17).Value = 2
$ws.Cells.Item(27, 18).Value = "Hortaliza"
